# Apply cryptos price/volume updates (commit: "Updated cryptos list on Mon May 27 05:07:09 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.771.66"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "3.917.64"
$ws.Range("E3").Value = "  +4.56%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'603.65"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "'165.77"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "3.918.15"
$ws.Range("E7").Value = "  +4.56%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.530"
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("E10").Value = "  -3.64%  "
$ws.Range("D11").Value = "'6.39"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "'0.462"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "'37.18"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").Value = "'0.0000245"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").Value = "4.574.46"
$ws.Range("E15").Value = "  +4.53%  "
$ws.Range("D16").Value = "3.950.16"
$ws.Range("E16").Value = "  +4.73%  "
$ws.Range("D17").Value = "68.934.35"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "'7.46"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").Value = "'17.01"
$ws.Range("E20").Value = "  -4.24%  "
$ws.Range("D21").Value = "'11.07"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").Value = "'486.96"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").Value = "'0.720"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("E24").Value = "  +11.93%  "
$ws.Range("D25").Value = "'84.39"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").Value = "'2.25"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").Value = "'12.08"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").Value = "'10.11"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("D31").Value = "4.071.56"
$ws.Range("E31").Value = "  +4.55%  "
$ws.Range("D32").Value = "'2.39"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("D34").Value = "'32.08"
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("D35").Value = "3.862.24"
$ws.Range("E35").Value = "  +4.82%  "
$ws.Range("D36").Value = "'0.107"
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").Value = "'5.89"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").Value = "'0.320"
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("D42").Value = "'3.00"
$ws.Range("E42").Value = "  -2.67%  "
$ws.Range("D43").Value = "'436.54"
$ws.Range("E43").Value = "  +2.99%  "
$ws.Range("D44").Value = "'1.99"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("D48").Value = "2.838.05"
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'142.13"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'26.10"
$ws.Range("E50").Value = "  +10.22%  "
$ws.Range("E51").Value = "  +1.20%  "
